$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43 (ALC)
$ws.Range("H43").Value = 834.4
$ws.Range("I43").Value = 544.3333
$ws.Range("J43").Value = 958.7143
$ws.Range("K43").Value = 544.3333
$ws.Range("L43").Value = 958.7143
$ws.Range("M43").Value = -475.3333
$ws.Range("N43").Value = -1096.7143

# Row 129 (ALC)
$ws.Range("H129").Value = 143600.4
$ws.Range("J129").Value = 154623.75
$ws.Range("L129").Value = 463871.25
$ws.Range("N129").Value = -473871.25

$ws = $wb.Worksheets.Item("ARM")
# Row 63 (ARM)
$ws.Range("H63").Value = 2606370
$ws.Range("I63").Value = 2394
$ws.Range("K63").Value = 2394
$ws.Range("M63").Value = -1708

# Row 66 (ARM)
$ws.Range("H66").Value = 2606370
$ws.Range("I66").Value = 2394
$ws.Range("K66").Value = 11970
$ws.Range("M66").Value = -8538

# Row 88 (ARM)
$ws.Range("H88").Value = 168275.67
$ws.Range("I88").Value = 933.3333
$ws.Range("J88").Value = 335618
$ws.Range("K88").Value = 933.3333
$ws.Range("L88").Value = 335618
$ws.Range("M88").Value = -527.3333
$ws.Range("N88").Value = -336430

# Row 91 (ARM)
$ws.Range("H91").Value = 168275.67
$ws.Range("I91").Value = 933.3333
$ws.Range("J91").Value = 335618
$ws.Range("K91").Value = 933.3333
$ws.Range("L91").Value = 335618
$ws.Range("M91").Value = 470.6667
$ws.Range("N91").Value = -338426

# Row 132 (ARM)
$ws.Range("H132").Value = 13564.688
$ws.Range("I132").Value = 2375.2424
$ws.Range("K132").Value = 7125.7272
$ws.Range("M132").Value = -4595.7272

$ws = $wb.Worksheets.Item("BSM")
# Row 35 (BSM)
$ws.Range("H35").Value = 24999.6
$ws.Range("J35").Value = 24999.6
$ws.Range("L35").Value = 24999.6
$ws.Range("N35").Value = -25619.6

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (CRP)
$ws.Range("H16").Value = 1174
$ws.Range("I16").Value = 982.1667
$ws.Range("J16").Value = 1749.5
$ws.Range("K16").Value = 982.1667
$ws.Range("L16").Value = 1749.5
$ws.Range("M16").Value = -695.1667
$ws.Range("N16").Value = -2323.5

# Row 99 (CRP)
$ws.Range("H99").Value = 3856.6667
$ws.Range("I99").Value = 2758.889
$ws.Range("K99").Value = 2758.889
$ws.Range("M99").Value = -1260.889

# Row 105 (CRP)
$ws.Range("H105").Value = 1089.2727
$ws.Range("I105").Value = 998.2
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 998.2
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 748.8
$ws.Range("N105").Value = -5494

# Row 113 (CRP)
$ws.Range("H113").Value = 1174
$ws.Range("I113").Value = 982.1667
$ws.Range("J113").Value = 1749.5
$ws.Range("K113").Value = 982.1667
$ws.Range("L113").Value = 1749.5
$ws.Range("M113").Value = 1187.8333
$ws.Range("N113").Value = -6089.5

# Row 122 (CRP)
$ws.Range("H122").Value = 1027.5
$ws.Range("I122").Value = 869.75
$ws.Range("J122").Value = 1216.8
$ws.Range("K122").Value = 2609.25
$ws.Range("L122").Value = 3650.4
$ws.Range("M122").Value = -159.25
$ws.Range("N122").Value = -8550.4

# Row 126 (CRP)
$ws.Range("H126").Value = 3856.6667
$ws.Range("I126").Value = 2758.889
$ws.Range("K126").Value = 8276.667000000001
$ws.Range("M126").Value = -5806.667000000001

# Row 132 (CRP)
$ws.Range("H132").Value = 3453.8572
$ws.Range("I132").Value = 2502.4285
$ws.Range("J132").Value = 5356.7144
$ws.Range("K132").Value = 7507.2855
$ws.Range("L132").Value = 16070.1432
$ws.Range("M132").Value = -4977.2855
$ws.Range("N132").Value = -21130.1432

# Row 134 (CRP)
$ws.Range("H134").Value = 761.9761999999999
$ws.Range("I134").Value = 703.1515000000001
$ws.Range("K134").Value = 2109.4545
$ws.Range("M134").Value = 425.5454999999997

$ws = $wb.Worksheets.Item("CUL")
# Row 2 (CUL)
$ws.Range("H2").Value = 9111.727999999999
$ws.Range("I2").Value = 16678.334
$ws.Range("J2").Value = 31.8
$ws.Range("K2").Value = 100070.004
$ws.Range("L2").Value = 190.8
$ws.Range("M2").Value = -99957.00399999999
$ws.Range("N2").Value = -416.8

# Row 23 (CUL)
$ws.Range("H23").Value = 331.33334
$ws.Range("J23").Value = 359.63635
$ws.Range("L23").Value = 1078.90905
$ws.Range("N23").Value = -1548.90905

# Row 68 (CUL)
$ws.Range("H68").Value = 1531.9333
$ws.Range("J68").Value = 1570
$ws.Range("L68").Value = 4710
$ws.Range("N68").Value = -6332

# Row 71 (CUL)
$ws.Range("H71").Value = 1531.9333
$ws.Range("J71").Value = 1570
$ws.Range("L71").Value = 14130
$ws.Range("N71").Value = -22242

# Row 131 (CUL)
$ws.Range("H131").Value = 722.6900000000001
$ws.Range("J131").Value = 740.0208
$ws.Range("L131").Value = 2220.0624
$ws.Range("N131").Value = -12300.0624

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (GSM)
$ws.Range("H80").Value = 3900
$ws.Range("I80").Value = 3840
$ws.Range("K80").Value = 3840
$ws.Range("M80").Value = -2842

# Row 83 (GSM)
$ws.Range("H83").Value = 3900
$ws.Range("I83").Value = 3840
$ws.Range("K83").Value = 19200
$ws.Range("M83").Value = -14208

# Row 122 (GSM)
$ws.Range("H122").Value = 5908.25
$ws.Range("I122").Value = 5549.875
$ws.Range("K122").Value = 16649.625
$ws.Range("M122").Value = -14199.625

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (LTW)
$ws.Range("H82").Value = 1621
$ws.Range("I82").Value = 1645.5555
$ws.Range("J82").Value = 1400
$ws.Range("K82").Value = 1645.5555
$ws.Range("L82").Value = 1400
$ws.Range("M82").Value = -1284.5555
$ws.Range("N82").Value = -2122

# Row 85 (LTW)
$ws.Range("H85").Value = 1621
$ws.Range("I85").Value = 1645.5555
$ws.Range("J85").Value = 1400
$ws.Range("K85").Value = 1645.5555
$ws.Range("L85").Value = 1400
$ws.Range("M85").Value = -397.5554999999999
$ws.Range("N85").Value = -3896

# Row 93 (LTW)
$ws.Range("H93").Value = 1249.1562
$ws.Range("I93").Value = 1136.5714
$ws.Range("K93").Value = 1136.5714
$ws.Range("M93").Value = 111.4286

# Row 132 (LTW)
$ws.Range("H132").Value = 319127.28
$ws.Range("J132").Value = 3598.75
$ws.Range("L132").Value = 10796.25
$ws.Range("N132").Value = -15856.25

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 2214.111
$ws.Range("J81").Value = 4306.5
$ws.Range("L81").Value = 8613
$ws.Range("N81").Value = -10735

# Row 84 (WVR)
$ws.Range("H84").Value = 2214.111
$ws.Range("J84").Value = 4306.5
$ws.Range("L84").Value = 43065
$ws.Range("N84").Value = -53673

# Row 113 (WVR)
$ws.Range("H113").Value = 1841.6154
$ws.Range("I113").Value = 1439.5416
$ws.Range("K113").Value = 4318.6248
$ws.Range("M113").Value = -2148.6248

# Row 122 (WVR)
$ws.Range("H122").Value = 1333.2667
$ws.Range("I122").Value = 1018.5455
$ws.Range("J122").Value = 2198.75
$ws.Range("K122").Value = 3055.6365
$ws.Range("L122").Value = 6596.25
$ws.Range("M122").Value = -605.6364999999996
$ws.Range("N122").Value = -11496.25

# Row 126 (WVR)
$ws.Range("H126").Value = 1982.4667
$ws.Range("I126").Value = 1457.25
$ws.Range("J126").Value = 4083.3333
$ws.Range("K126").Value = 4371.75
$ws.Range("L126").Value = 12249.9999
$ws.Range("M126").Value = -1901.75
$ws.Range("N126").Value = -17189.9999
